# Updated cryptos list on Sat Oct 28 10:33:50 UTC 2023 with GitHub Actions
#
# Price (col D) / Volume(1h) (col E) refresh for the coinranking.com scrape.
# A couple of coins (WrappedBTC/Polygon and Kaspa/FraxShare) swapped rank
# positions between rows 15/16 and 45/46, so those rows' Coin name, Link and
# Price/Volume cells are all rewritten to their new row.
#
# Note: several new Price values are plain decimal numbers (e.g. "227.55",
# "10.90") that Excel would otherwise auto-convert to floating point,
# silently dropping the trailing zero / exact text form. A leading
# apostrophe forces those assignments to be stored as text, matching the
# original cell's text semantics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.102.48'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '1.790.27'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''227.55'
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '''32.29'
$ws.Range('E8').Value = '  -1.36%  '
$ws.Range('E9').Value = '  +2.91%  '
$ws.Range('D10').Value = '''0.0691'
$ws.Range('E10').Value = '  -2.80%  '
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').Value = '2.047.84'
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '''11.55'
$ws.Range('E13').Value = '  +5.05%  '
$ws.Range('D14').Value = '1.791.84'
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '''0.622'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '34.090.39'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('E17').Value = '  +0.72%  '
$ws.Range('D18').Value = '''67.85'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '''245.05'
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('D22').Value = '''10.90'
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('D23').Value = '''4.11'
$ws.Range('E23').Value = '  +0.59%  '
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('D25').Value = '''161.99'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').Value = '''7.17'
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('D27').Value = '''16.32'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  +1.09%  '
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('E30').Value = '  +2.31%  '
$ws.Range('E31').Value = '  +1.70%  '
$ws.Range('D32').Value = '''3.67'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').Value = '''3.63'
$ws.Range('E33').Value = '  +3.26%  '
$ws.Range('E34').Value = '  +0.76%  '
$ws.Range('D35').Value = '1.423.73'
$ws.Range('E35').Value = '  +2.16%  '
$ws.Range('D36').Value = '''0.644'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +2.62%  '
$ws.Range('E38').Value = '  +7.31%  '
$ws.Range('E39').Value = '  -0.62%  '
$ws.Range('D40').Value = '''80.62'
$ws.Range('E40').Value = '  +3.18%  '
$ws.Range('E41').Value = '  +0.44%  '
$ws.Range('D42').Value = '''0.922'
$ws.Range('E42').Value = '  +0.79%  '
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('D44').Value = '''13.39'
$ws.Range('E44').Value = '  +8.00%  '
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').Value = '''0.0509'
$ws.Range('E45').Value = '  +2.33%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = '''6.06'
$ws.Range('E46').Value = '  +3.74%  '
$ws.Range('D47').Value = '0.0₆0138'
$ws.Range('E47').Value = '  -1.59%  '
$ws.Range('E48').Value = '  -1.06%  '
$ws.Range('D49').Value = '''107.59'
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').Value = '1.948.76'
$ws.Range('E50').Value = '  +0.21%  '
